$wb = $excel.ActiveWorkbook

# Sheet ALC, row 69
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 7062.5
$ws.Range("I69").Value = 50000
$ws.Range("J69").Value = 4200
$ws.Range("K69").Value = 150000
$ws.Range("L69").Value = 12600
$ws.Range("M69").Value = -149126
$ws.Range("N69").Value = -14348

# Sheet ALC, row 72
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 7062.5
$ws.Range("I72").Value = 50000
$ws.Range("J72").Value = 4200
$ws.Range("K72").Value = 450000
$ws.Range("L72").Value = 37800
$ws.Range("M72").Value = -445632
$ws.Range("N72").Value = -46536

# Sheet ALC, row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 43481172
$ws.Range("I76").Value = 55558388
$ws.Range("J76").Value = 3200
$ws.Range("K76").Value = 55558388
$ws.Range("L76").Value = 3200
$ws.Range("M76").Value = -55558073
$ws.Range("N76").Value = -3830

# Sheet ALC, row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 43481172
$ws.Range("I79").Value = 55558388
$ws.Range("J79").Value = 3200
$ws.Range("K79").Value = 55558388
$ws.Range("L79").Value = 3200
$ws.Range("M79").Value = -55557296
$ws.Range("N79").Value = -5384

# Sheet ALC, row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 3576743
$ws.Range("I80").Value = 9587.462
$ws.Range("J80").Value = 5684608
$ws.Range("K80").Value = 28762.386
$ws.Range("L80").Value = 17053824
$ws.Range("M80").Value = -27764.386
$ws.Range("N80").Value = -17055820

# Sheet ALC, row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 3576743
$ws.Range("I83").Value = 9587.462
$ws.Range("J83").Value = 5684608
$ws.Range("K83").Value = 86287.158
$ws.Range("L83").Value = 51161472
$ws.Range("M83").Value = -81295.158
$ws.Range("N83").Value = -51171456

# Sheet ALC, row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 10102593
$ws.Range("I88").Value = 1001.5
$ws.Range("J88").Value = 12988763
$ws.Range("K88").Value = 1001.5
$ws.Range("L88").Value = 12988763
$ws.Range("M88").Value = -595.5
$ws.Range("N88").Value = -12989575

# Sheet ALC, row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 10102593
$ws.Range("I91").Value = 1001.5
$ws.Range("J91").Value = 12988763
$ws.Range("K91").Value = 1001.5
$ws.Range("L91").Value = 12988763
$ws.Range("M91").Value = 402.5
$ws.Range("N91").Value = -12991571

# Sheet ALC, row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 346679.97
$ws.Range("I129").Value = 494.2353
$ws.Range("J129").Value = 837109.75
$ws.Range("K129").Value = 1482.7059
$ws.Range("L129").Value = 2511329.25
$ws.Range("M129").Value = 3517.2941
$ws.Range("N129").Value = -2521329.25

# Sheet ALC, row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2268886.8
$ws.Range("I132").Value = 1340.3125
$ws.Range("J132").Value = 111111110
$ws.Range("K132").Value = 4020.9375
$ws.Range("L132").Value = 333333330
$ws.Range("M132").Value = -1490.9375
$ws.Range("N132").Value = -333338390

# Sheet ARM, row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 88350
$ws.Range("I2").Value = 117200
$ws.Range("K2").Value = 117200
$ws.Range("M2").Value = -117087

# Sheet ARM, row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 401043.25
$ws.Range("I45").Value = 910017.4
$ws.Range("J45").Value = 1135
$ws.Range("K45").Value = 910017.4
$ws.Range("L45").Value = 1135
$ws.Range("M45").Value = -909640.4
$ws.Range("N45").Value = -1889

# Sheet ARM, row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 526.9375
$ws.Range("I97").Value = 455.3846
$ws.Range("J97").Value = 837
$ws.Range("K97").Value = 455.3846
$ws.Range("L97").Value = 837
$ws.Range("M97").Value = 40.61540000000002
$ws.Range("N97").Value = -1829

# Sheet ARM, row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 88350
$ws.Range("I116").Value = 117200
$ws.Range("K116").Value = 117200
$ws.Range("M116").Value = -114906

# Sheet BSM, row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 88350
$ws.Range("I3").Value = 117200
$ws.Range("K3").Value = 117200
$ws.Range("M3").Value = -117086

# Sheet BSM, row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1970.98
$ws.Range("I86").Value = 1980.5918
$ws.Range("J86").Value = 1500
$ws.Range("K86").Value = 1980.5918
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = -857.5917999999999
$ws.Range("N86").Value = -3746

# Sheet BSM, row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1970.98
$ws.Range("I89").Value = 1980.5918
$ws.Range("J89").Value = 1500
$ws.Range("K89").Value = 9902.958999999999
$ws.Range("L89").Value = 7500
$ws.Range("M89").Value = -4286.958999999999
$ws.Range("N89").Value = -18732

# Sheet CRP, row 130
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H130").Value = 46913.332
$ws.Range("J130").Value = 46913.332
$ws.Range("L130").Value = 46913.332
$ws.Range("N130").Value = -56953.332

# Sheet CRP, row 131
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H131").Value = 16098.667
$ws.Range("J131").Value = 16098.667
$ws.Range("L131").Value = 16098.667
$ws.Range("N131").Value = -26178.667

# Sheet CUL, row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1808187.2
$ws.Range("I5").Value = 1972945.6
$ws.Range("J5").Value = 1516115.8
$ws.Range("K5").Value = 5918836.800000001
$ws.Range("L5").Value = 4548347.4
$ws.Range("M5").Value = -5918724.800000001
$ws.Range("N5").Value = -4548571.4

# Sheet CUL, row 14
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 525.4583
$ws.Range("I14").Value = 525.4583
$ws.Range("K14").Value = 1576.3749
$ws.Range("M14").Value = -1403.3749

# Sheet CUL, row 117
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1504.0952
$ws.Range("I117").Value = 729
$ws.Range("J117").Value = 1542.85
$ws.Range("K117").Value = 2187
$ws.Range("L117").Value = 4628.549999999999
$ws.Range("M117").Value = 1255
$ws.Range("N117").Value = -11512.55

# Sheet CUL, row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1312.16
$ws.Range("I122").Value = 278.82352
$ws.Range("K122").Value = 2509.41168
$ws.Range("M122").Value = -59.41167999999971

# Sheet CUL, row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2635.875
$ws.Range("I129").Value = 2442.3076
$ws.Range("J129").Value = 2768.3157
$ws.Range("K129").Value = 7326.9228
$ws.Range("L129").Value = 8304.947100000001
$ws.Range("M129").Value = -2326.9228
$ws.Range("N129").Value = -18304.9471

# Sheet CUL, row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1808187.2
$ws.Range("I135").Value = 1972945.6
$ws.Range("J135").Value = 1516115.8
$ws.Range("K135").Value = 17756510.4
$ws.Range("L135").Value = 13645042.2
$ws.Range("M135").Value = -17753975.4
$ws.Range("N135").Value = -13650112.2

# Sheet GSM, row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 241.11111
$ws.Range("I107").Value = 130
$ws.Range("J107").Value = 463.33334
$ws.Range("K107").Value = 130
$ws.Range("L107").Value = 463.33334
$ws.Range("M107").Value = 1790
$ws.Range("N107").Value = -4303.33334

# Sheet LTW, row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1599.0476
$ws.Range("I61").Value = 1056.1538
$ws.Range("J61").Value = 2481.25
$ws.Range("K61").Value = 1056.1538
$ws.Range("L61").Value = 2481.25
$ws.Range("M61").Value = -854.1538
$ws.Range("N61").Value = -2885.25

# Sheet LTW, row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 8218.703
$ws.Range("I93").Value = 2289.9
$ws.Range("K93").Value = 2289.9
$ws.Range("M93").Value = -1041.9

# Sheet LTW, row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1599.0476
$ws.Range("I113").Value = 1056.1538
$ws.Range("J113").Value = 2481.25
$ws.Range("K113").Value = 1056.1538
$ws.Range("L113").Value = 2481.25
$ws.Range("M113").Value = 1113.8462
$ws.Range("N113").Value = -6821.25

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3183594
$ws.Range("I132").Value = 5565806.5
$ws.Range("J132").Value = 7310.8887
$ws.Range("K132").Value = 16697419.5
$ws.Range("L132").Value = 21932.6661
$ws.Range("M132").Value = -16694889.5
$ws.Range("N132").Value = -26992.6661

# Sheet WVR, row 54
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 99166
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 99166
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 99166
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -100206

# Sheet WVR, row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 20929.5
$ws.Range("I81").Value = 1065.6666
$ws.Range("J81").Value = 26346.908
$ws.Range("K81").Value = 2131.3332
$ws.Range("L81").Value = 52693.816
$ws.Range("M81").Value = -1070.3332
$ws.Range("N81").Value = -54815.816

# Sheet WVR, row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 20929.5
$ws.Range("I84").Value = 1065.6666
$ws.Range("J84").Value = 26346.908
$ws.Range("K84").Value = 10656.666
$ws.Range("L84").Value = 263469.08
$ws.Range("M84").Value = -5352.666000000001
$ws.Range("N84").Value = -274077.08

# Sheet WVR, row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 462.5
$ws.Range("J113").Value = 505
$ws.Range("L113").Value = 1515
$ws.Range("N113").Value = -5855

# Sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1776.1177
$ws.Range("I122").Value = 1638.4
$ws.Range("J122").Value = 1972.8572
$ws.Range("K122").Value = 4915.200000000001
$ws.Range("L122").Value = 5918.571599999999
$ws.Range("M122").Value = -2465.200000000001
$ws.Range("N122").Value = -10818.5716
